# Update the account-statement worksheet: swap the order of the "Periodo Mora"
# values shown in E16:E19 (1704,1705,1711,1712 -> 1712,1711,1705,1704) and
# keep the "Valor Mora" (F column) figures attached to the correct period.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E16").Value = "1712"
$ws.Range("F16").Value = 11808

$ws.Range("E17").Value = "1711"
$ws.Range("F17").Value = 29520

$ws.Range("E18").Value = "1705"
$ws.Range("F18").Value = 29520

$ws.Range("E19").Value = "1704"
$ws.Range("F19").Value = 29520
